# Apply a cyclic re-ordering of species-observation data across rows 9-15.
# Only columns A, B, E, F, G, H, Q, R (Id, Taxonsorteringsordning, TaxonId,
# Artnamn, Vetenskapligt namn, Auktor, Ost, Nord) move between rows; all the
# other columns stay tied to their original row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that are rearranged between the affected rows.
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

# Rows participating in the shuffle.
$rows = @(9, 10, 11, 12, 13, 14, 15)

# Snapshot the current ("before") values for each affected row/column so we
# can redistribute them without relying on hard-coded literals.
$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Destination row -> source row mapping (i.e. new row N receives the data
# that used to live in row mapping[N]).
$mapping = @{
    9  = 12
    10 = 13
    11 = 14
    12 = 9
    13 = 15
    14 = 10
    15 = 11
}

foreach ($destRow in $rows) {
    $srcRow = $mapping[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $srcData[$col]
    }
}
